$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "System, *") {
        $parts = $val -split ", "
        $newParts = $parts[1..($parts.Length - 1)] + @($parts[0])
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value2 = $newVal
    }
}
